# Applies the E2EComponent.pptx edit:
#  1. Refresh the cached "datetimeFigureOut" field text (footer date stamps)
#     on the Notes Master and on every Slide Layout's Date placeholder -
#     the deck was re-saved on 14/8/2020 (was 28/9/2018 / 9/28/2018).
#  2. Rename "::scalability" -> "::lnp" (split into two runs, same as
#     PowerPoint does when the trailing word is retyped) on the small
#     package-path label shape.
#  3. Rename the "Scalability tests" box to "L&P tests".

$p = $ppt.ActivePresentation

# --- 1a. Notes Master date placeholder (dd/m/yyyy style: 28/9/2018 -> 14/8/2020)
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1 -and $shp.TextFrame.TextRange.Text -eq "28/9/2018") {
        $shp.TextFrame.TextRange.Text = "14/8/2020"
    }
}

# --- 1b. Every Slide Layout's date placeholder (m/d/yyyy style: 9/28/2018 -> 8/14/2020)
$master = $p.SlideMaster
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame -eq -1 -and $shp.TextFrame.TextRange.Text -eq "9/28/2018") {
            $shp.TextFrame.TextRange.Text = "8/14/2020"
        }
    }
}

# --- 2. "::scalability" -> "::lnp" (two runs: "::" then "lnp")
$slide = $p.Slides.Item(1)
$pathGroup = $slide.Shapes.Item(7)
$labelShape = $pathGroup.GroupItems.Item(1)
$labelShape.TextFrame.TextRange.Text = "::lnp"
$labelShape.TextFrame.TextRange.Characters(3, 3).Font.Bold = -1

# --- 3. "Scalability tests" -> "L&P tests"
$testsShape = $slide.Shapes.Item(8)
$testsShape.TextFrame.TextRange.Text = "L&P tests"
